$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: single value change ---
$ws.Range("B14").Value = 89072

# --- Rows 15 and 16 have their entire contents swapped ---

# New row 15 values (previously on row 16)
$ws.Range("A15").Value = 112481246
$ws.Range("B15").Value = 90291
$ws.Range("C15").Value = "Ovaliderad"
$ws.Range("D15").Value = "VU"
$ws.Range("E15").Value = 1958
$ws.Range("F15").Value = "Lammticka"
$ws.Range("G15").Value = "Albatrellus subrubescens"
$ws.Range("H15").Value = "(Murrill) Pouzar"
$ws.Range("I15").ClearContents()
$ws.Range("J15").ClearContents()
$ws.Range("K15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("P15").Value = "Barrtjärnen, Bergsjövägen, Ede, Nordanstigs kommun (Barrtjärnen, Bergsjövägen, Ede, Nordanstigs kommun), Hls"
$ws.Range("Q15").Value = 601210
$ws.Range("R15").Value = 6877659
$ws.Range("S15").Value = 25
$ws.Range("T15").Value = "Gävleborg"
$ws.Range("U15").Value = "Nordanstig"
$ws.Range("V15").Value = "Hälsingland"
$ws.Range("W15").Value = "Bergsjö"
$ws.Range("Y15").NumberFormat = "@"
$ws.Range("Y15").Value = "2023-10-02"
$ws.Range("Z15").Value = "16:05"
$ws.Range("AA15").NumberFormat = "@"
$ws.Range("AA15").Value = "2023-10-02"
$ws.Range("AB15").Value = "16:05"
$ws.Range("AC15").Value = "I kanten av stigen ned mot båtplats. Västra kanten I böjen."
$ws.Range("AD15").Value = $false
$ws.Range("AE15").Value = $false
$ws.Range("AF15").ClearContents()
$ws.Range("AG15").Value = $false
$ws.Range("AT15").ClearContents()
$ws.Range("AW15").Value = "Henrik Tykosson"
$ws.Range("AX15").Value = "Henrik Tykosson"
$ws.Range("AY15").ClearContents()

# New row 16 values (previously on row 15)
$ws.Range("A16").Value = 112481511
$ws.Range("B16").Value = 90813
$ws.Range("C16").Value = "Ovaliderad"
$ws.Range("D16").Value = "VU"
$ws.Range("E16").Value = 1435
$ws.Range("F16").Value = "Bitter taggsvamp"
$ws.Range("G16").Value = "Hydnellum fennicum"
$ws.Range("H16").Value = "(P.Karst.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "1"
$ws.Range("J16").Value = "fruktkroppar"
$ws.Range("K16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("P16").Value = "Barrtjärnen, Bergsjövägen, Ede, Nordanstigs kommun (Barrtjärnen, Bergsjövägen, Ede, Nordanstigs kommun), Hls"
$ws.Range("Q16").Value = 601183
$ws.Range("R16").Value = 6877672
$ws.Range("S16").Value = 25
$ws.Range("T16").Value = "Gävleborg"
$ws.Range("U16").Value = "Nordanstig"
$ws.Range("V16").Value = "Hälsingland"
$ws.Range("W16").Value = "Bergsjö"
$ws.Range("Y16").NumberFormat = "@"
$ws.Range("Y16").Value = "2023-09-09"
$ws.Range("Z16").Value = "15:00"
$ws.Range("AA16").NumberFormat = "@"
$ws.Range("AA16").Value = "2023-09-09"
$ws.Range("AB16").Value = "15:00"
$ws.Range("AC16").Value = "I slänten ner mot en större svacka i terrängen. Ca 3 m från stigen. Tydlig doft av bittermandel."
$ws.Range("AD16").Value = $false
$ws.Range("AE16").Value = $false
$ws.Range("AF16").ClearContents()
$ws.Range("AG16").Value = $false
$ws.Range("AT16").ClearContents()
$ws.Range("AW16").Value = "Henrik Tykosson"
$ws.Range("AX16").Value = "Henrik Tykosson"
$ws.Range("AY16").ClearContents()
